$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 93
$ws.Range("I11").Value = 93
$ws.Range("K11").Value = 93
$ws.Range("M11").Value = 47
$ws.Range("H87").Value = 49997.5
$ws.Range("J87").Value = 49997.5
$ws.Range("L87").Value = 49997.5
$ws.Range("N87").Value = -52493.5
$ws.Range("H90").Value = 49997.5
$ws.Range("J90").Value = 49997.5
$ws.Range("L90").Value = 149992.5
$ws.Range("N90").Value = -162472.5
$ws.Range("H131").Value = 209.2
$ws.Range("I131").Value = 209.2
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 627.5999999999999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 4412.4
$ws.Range("N131").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 4024.75
$ws.Range("I19").Value = 2033
$ws.Range("K19").Value = 2033
$ws.Range("M19").Value = -1804
$ws.Range("H32").Value = 3868.907
$ws.Range("I32").Value = 1445.4103
$ws.Range("J32").Value = 27498
$ws.Range("K32").Value = 1445.4103
$ws.Range("L32").Value = 27498
$ws.Range("M32").Value = -1158.4103
$ws.Range("N32").Value = -28072
$ws.Range("H74").Value = 2427.7856
$ws.Range("J74").Value = 4999.5
$ws.Range("L74").Value = 4999.5
$ws.Range("N74").Value = -6747.5
$ws.Range("H77").Value = 2427.7856
$ws.Range("J77").Value = 4999.5
$ws.Range("L77").Value = 24997.5
$ws.Range("N77").Value = -33733.5
$ws.Range("H132").Value = 2193
$ws.Range("I132").Value = 2193
$ws.Range("K132").Value = 6579
$ws.Range("M132").Value = -4049

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 29976.666
$ws.Range("I15").Value = 18000
$ws.Range("J15").Value = 35965
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 35965
$ws.Range("M15").Value = -17773
$ws.Range("N15").Value = -36419
$ws.Range("H19").Value = 25567.834
$ws.Range("I19").Value = 25801.273
$ws.Range("J19").Value = 23000
$ws.Range("K19").Value = 25801.273
$ws.Range("L19").Value = 23000
$ws.Range("M19").Value = -25628.273
$ws.Range("N19").Value = -23346
$ws.Range("H20").Value = 9002.333000000001
$ws.Range("I20").Value = 9002.333000000001
$ws.Range("K20").Value = 9002.333000000001
$ws.Range("M20").Value = -8755.333000000001
$ws.Range("H56").Value = 100100
$ws.Range("I56").Value = 100100
$ws.Range("K56").Value = 100100
$ws.Range("M56").Value = -99361
$ws.Range("H75").Value = 71000
$ws.Range("J75").Value = 130000
$ws.Range("L75").Value = 130000
$ws.Range("N75").Value = -131872
$ws.Range("H78").Value = 71000
$ws.Range("J78").Value = 130000
$ws.Range("L78").Value = 390000
$ws.Range("N78").Value = -399360
$ws.Range("H86").Value = 4356.8335
$ws.Range("I86").Value = 4033.7778
$ws.Range("J86").Value = 5326
$ws.Range("K86").Value = 4033.7778
$ws.Range("L86").Value = 5326
$ws.Range("M86").Value = -2910.7778
$ws.Range("N86").Value = -7572
$ws.Range("H89").Value = 4356.8335
$ws.Range("I89").Value = 4033.7778
$ws.Range("J89").Value = 5326
$ws.Range("K89").Value = 20168.889
$ws.Range("L89").Value = 26630
$ws.Range("M89").Value = -14552.889
$ws.Range("N89").Value = -37862

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 3500
$ws.Range("J8").Value = 3500
$ws.Range("L8").Value = 3500
$ws.Range("N8").Value = -3780
$ws.Range("H14").Value = 8677.875
$ws.Range("I14").Value = 439.5
$ws.Range("J14").Value = 11424
$ws.Range("K14").Value = 439.5
$ws.Range("L14").Value = 11424
$ws.Range("M14").Value = -269.5
$ws.Range("N14").Value = -11764
$ws.Range("H19").Value = 3268.5625
$ws.Range("I19").Value = 1319
$ws.Range("J19").Value = 5218.125
$ws.Range("K19").Value = 1319
$ws.Range("L19").Value = 5218.125
$ws.Range("M19").Value = -1149
$ws.Range("N19").Value = -5558.125
$ws.Range("H24").Value = 3268.5625
$ws.Range("I24").Value = 1319
$ws.Range("J24").Value = 5218.125
$ws.Range("K24").Value = 1319
$ws.Range("L24").Value = 5218.125
$ws.Range("M24").Value = -1149
$ws.Range("N24").Value = -5558.125
$ws.Range("H45").Value = 52000
$ws.Range("J45").Value = 52000
$ws.Range("L45").Value = 52000
$ws.Range("N45").Value = -53186
$ws.Range("H50").Value = 20061.75
$ws.Range("J50").Value = 19996.875
$ws.Range("L50").Value = 19996.875
$ws.Range("N50").Value = -21246.875
$ws.Range("H59").Value = 29214.715
$ws.Range("J59").Value = 29999.834
$ws.Range("L59").Value = 29999.834
$ws.Range("N59").Value = -32289.834
$ws.Range("H60").Value = 21206.7
$ws.Range("J60").Value = 19996.428
$ws.Range("L60").Value = 19996.428
$ws.Range("N60").Value = -21018.428
$ws.Range("H68").Value = 27499.875
$ws.Range("H71").Value = 27499.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4294.758
$ws.Range("I4").Value = 5509.8696
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 16529.6088
$ws.Range("L4").Value = 4500
$ws.Range("M4").Value = -16417.6088
$ws.Range("N4").Value = -4724
$ws.Range("H134").Value = 169065.5
$ws.Range("I134").Value = 201278.6
$ws.Range("K134").Value = 603835.8
$ws.Range("M134").Value = -598765.8
$ws.Range("H136").Value = 2880.6
$ws.Range("I136").Value = 2880.6
$ws.Range("K136").Value = 8641.799999999999
$ws.Range("M136").Value = -3541.799999999999
$ws.Range("H137").Value = 2000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null
$ws.Range("H138").Value = 2999.8333
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280
$ws.Range("H139").Value = 1343.6364
$ws.Range("I139").Value = 878.3
$ws.Range("K139").Value = 2634.9
$ws.Range("M139").Value = 2505.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3976.5293
$ws.Range("I43").Value = 2621.9
$ws.Range("J43").Value = 5911.7144
$ws.Range("K43").Value = 2621.9
$ws.Range("L43").Value = 5911.7144
$ws.Range("M43").Value = -2470.9
$ws.Range("N43").Value = -6213.7144
$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21640
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = $null
$ws.Range("H126").Value = 11112.25
$ws.Range("I126").Value = 6474.5
$ws.Range("J126").Value = 15750
$ws.Range("K126").Value = 19423.5
$ws.Range("L126").Value = 47250
$ws.Range("M126").Value = -16953.5
$ws.Range("N126").Value = -52190
$ws.Range("H133").Value = 39999.2
$ws.Range("J133").Value = 39999.2
$ws.Range("L133").Value = 39999.2
$ws.Range("N133").Value = -50119.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22589.295
$ws.Range("I7").Value = 22503.908
$ws.Range("J7").Value = 22745.834
$ws.Range("K7").Value = 22503.908
$ws.Range("L7").Value = 22745.834
$ws.Range("M7").Value = -22391.908
$ws.Range("N7").Value = -22969.834
$ws.Range("H126").Value = 22589.295
$ws.Range("I126").Value = 22503.908
$ws.Range("J126").Value = 22745.834
$ws.Range("K126").Value = 67511.724
$ws.Range("L126").Value = 68237.50199999999
$ws.Range("M126").Value = -65041.724
$ws.Range("N126").Value = -73177.50199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 2511247.8
$ws.Range("J15").Value = 14997.333
$ws.Range("L15").Value = 14997.333
$ws.Range("N15").Value = -15573.333
$ws.Range("H54").Value = 616000
$ws.Range("J54").Value = 765000
$ws.Range("L54").Value = 765000
$ws.Range("N54").Value = -766040
$ws.Range("H126").Value = 21134.393
$ws.Range("I126").Value = 22139.6
$ws.Range("J126").Value = 18621.375
$ws.Range("K126").Value = 66418.79999999999
$ws.Range("L126").Value = 55864.125
$ws.Range("M126").Value = -63948.79999999999
$ws.Range("N126").Value = -60804.125
$ws.Range("H132").Value = 2326.3572
$ws.Range("I132").Value = 811.625
$ws.Range("K132").Value = 2434.875
$ws.Range("M132").Value = 95.125
